$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume refresh. Column D ("Price") values that parse as plain
# numbers need the cell forced to Text format first, otherwise Excel would
# silently convert strings like "218.03" or "0.05150" into numbers and drop
# the formatting-significant trailing digits.

$ws.Range('D2').Value = '26.125.56'
$ws.Range('E2').Value = '  +0.61%  '

$ws.Range('D3').Value = '1.655.85'
$ws.Range('E3').Value = '  +0.55%  '

$ws.Range('E4').Value = '  -0.05%  '

$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '218.03'
$ws.Range('E5').Value = '  +0.61%  '

$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '0.5308'
$ws.Range('E6').Value = '  +2.53%  '

$ws.Range('E7').Value = '  -0.07%  '

$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.2616'
$ws.Range('E8').Value = '  +0.11%  '

$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.06333'
$ws.Range('E9').Value = '  +1.49%  '

$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '20.45'
$ws.Range('E10').Value = '  -0.06%  '

$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '0.07795'
$ws.Range('E11').Value = '  +1.13%  '

$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '4.523'
$ws.Range('E12').Value = '  +1.72%  '

$ws.Range('D13').Value = '1.656.55'
$ws.Range('E13').Value = '  +0.25%  '

$ws.Range('D14').Value = '1.882.78'
$ws.Range('E14').Value = '  +0.44%  '

$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.5496'
$ws.Range('E15').Value = '  +1.74%  '

$ws.Range('D16').Value = '0.0₅8220'
$ws.Range('E16').Value = '  +1.84%  '

$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '65.45'
$ws.Range('E17').Value = '  +1.38%  '

$ws.Range('D18').Value = '26.118.58'
$ws.Range('E18').Value = '  +0.64%  '

$ws.Range('E19').Value = '  -0.01%  '

$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '4.602'
$ws.Range('E20').Value = '  +1.13%  '

$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '191.32'
$ws.Range('E21').Value = '  +0.46%  '

$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '10.09'
$ws.Range('E22').Value = '  +1.23%  '

$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '6.034'
$ws.Range('E23').Value = '  +1.16%  '

$ws.Range('E24').Value = '  -0.07%  '

$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '145.09'
$ws.Range('E25').Value = '  +5.06%  '

$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '0.1228'
$ws.Range('E26').Value = '  +0.10%  '

$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '7.223'
$ws.Range('E27').Value = '  +0.24%  '

$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '16.01'
$ws.Range('E28').Value = '  +0.04%  '

$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.459'
$ws.Range('E29').Value = '  +4.01%  '

$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '0.05792'
$ws.Range('E30').Value = '  -2.01%  '

$ws.Range('E31').Value = '  +0.11%  '

$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '3.562'
$ws.Range('E32').Value = '  +1.59%  '

$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '3.278'
$ws.Range('E33').Value = '  +1.23%  '

$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '1.603'
$ws.Range('E34').Value = '  +3.45%  '

$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '2.803'
$ws.Range('E35').Value = '  +2.11%  '

$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '0.9524'
$ws.Range('E36').Value = '  +0.68%  '

$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '2.415'
$ws.Range('E37').Value = '  -0.11%  '

$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.5765'
$ws.Range('E38').Value = '  +2.45%  '

$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.01612'
$ws.Range('E39').Value = '  +1.62%  '

$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.8557'
$ws.Range('E40').Value = '  +1.41%  '

$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '5.805'
$ws.Range('E41').Value = '  -1.03%  '

$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '104.80'
$ws.Range('E42').Value = '  +4.23%  '

$ws.Range('B43').Value = 'Maker'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D43').Value = '1.046.68'
$ws.Range('E43').Value = '  +5.15%  '

$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '1.003'
$ws.Range('E44').Value = '  -0.08%  '

$ws.Range('D45').Value = '1.796.82'
$ws.Range('E45').Value = '  +0.33%  '

$ws.Range('E46').Value = '  +1.48%  '

$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '1.006'
$ws.Range('E47').Value = '  -0.24%  '

$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '0.4334'
$ws.Range('E48').Value = '  +0.74%  '

$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '7.885'
$ws.Range('E49').Value = '  -0.91%  '

$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.05150'
$ws.Range('E50').Value = '  +0.15%  '

$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '1.450'
$ws.Range('E51').Value = '  -0.97%  '
